$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose text values would otherwise be auto-coerced by Excel
# (purely numeric-looking text, or date-looking text) must be forced
# to Text format *before* assigning the value so they stay strings.
$ws.Range("I34").NumberFormat = "@"
$ws.Range("Y34").NumberFormat = "@"
$ws.Range("AA34").NumberFormat = "@"

# Numeric columns
$ws.Range("A34").Value = 112269372
$ws.Range("B34").Value = 90844
$ws.Range("E34").Value = 5449
$ws.Range("Q34").Value = 447837
$ws.Range("R34").Value = 6430211
$ws.Range("S34").Value = 25

# Text columns
$ws.Range("C34").Value = "Ovaliderad"
$ws.Range("D34").Value = "NT"
$ws.Range("F34").Value = "Svart taggsvamp"
$ws.Range("G34").Value = "Phellodon niger"
$ws.Range("H34").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("I34").Value = "3"
$ws.Range("J34").Value = "fruktkroppar"
$ws.Range("P34").Value = "800 m V Axtorp, Vg"
$ws.Range("T34").Value = "Jönköping"
$ws.Range("U34").Value = "Habo"
$ws.Range("V34").Value = "Västergötland"
$ws.Range("W34").Value = "Gustav Adolf"
$ws.Range("Y34").Value = "2023-09-22"
$ws.Range("AA34").Value = "2023-09-22"
$ws.Range("AW34").Value = "Niklas Johansson"
$ws.Range("AX34").Value = "Niklas Johansson"

# Boolean columns
$ws.Range("AD34").Value = $false
$ws.Range("AE34").Value = $false
$ws.Range("AG34").Value = $false

Write-Output "Row 34 added"
